$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91..134 down to 92..135.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new weekly record. It mirrors the
# (now shifted-down) row 92 except for the date and the volume/price columns.
$ws.Range("A91").Value = 6
$ws.Range("B91").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C91").Value = "Metropolitana"
$ws.Range("D91").Value = 44960
$ws.Range("E91").Value = 13
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100101
$ws.Range("H91").Value = "Berries"
$ws.Range("I91").Value = 100101008
$ws.Range("J91").Value = "Mora"
$ws.Range("K91").Value = "Sin especificar"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 175
$ws.Range("N91").Value = 3000
$ws.Range("O91").Value = 3000
$ws.Range("P91").Value = 3000
$ws.Range("Q91").Value = "$/bandeja 2 kilos"
$ws.Range("R91").Value = "Provincia de Curicó"
$ws.Range("S91").Value = 1500
$ws.Range("T91").Value = 2
